$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "name=Item_Code,dataType=text,title=Code,updateCriteria=true"
$ws.Range("G2").Value = "name=Protection,dataType=number,title=Protection"
$ws.Range("H2").Value = "name=Item_Value,dataType=number,title=Value"
$ws.Range("I2").Value = "name=Type,dataType=entity:Item_Types.Type"
$ws.Range("J2").Value = "name=Tech_Age,dataType=entity:Technological_Ages.Age"
$ws.Range("K2").Value = "name=Released,dataType=checkbox,title=Released"
$ws.Range("L2").Value = "name=Weapon,dataType=checkbox,title=Weapon"
$ws.Range("M2").Value = "name=Armor,dataType=checkbox,title=Armor"
$ws.Range("N2").Value = "name=Shield,dataType=checkbox,title=Shield"
